$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 22:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 952995
$ws.Range("C4").Value = 27763
$ws.Range("D4").Value = 116015
$ws.Range("E4").Value = 783235
$ws.Range("F4").Value = 15110
$ws.Range("G4").Value = 1552
$ws.Range("H4").Value = 53745

# Row 14 - Brasil
$ws.Range("B14").Value = 57961
$ws.Range("C14").Value = 4966
$ws.Range("E14").Value = 26343
$ws.Range("G14").Value = 293
$ws.Range("H14").Value = 3963

# Row 16 - Canada
$ws.Range("D16").Value = 16320
$ws.Range("E16").Value = 26234

# Row 99 - Niger
$ws.Range("B99").Value = 684
$ws.Range("C99").Value = 3
$ws.Range("D99").Value = 325
$ws.Range("E99").Value = 332
$ws.Range("G99").Value = 3
$ws.Range("H99").Value = 27
